$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column A, shifting existing data to column C
$ws.Range("A:B").EntireColumn.Insert()

# Set header row values
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

# Set data row values
$ws.Range("A2").Value = "abc"
$ws.Range("B2").Value = "xyz"

# Apply the same style as C1 (yellow fill) to the new header cells
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E10").Select()
